$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 'Test', 'Status', 'Timestamp'),
    @(2, 'Button Clicked', 'Pass', '2025-12-23T04:56:21.237Z'),
    @(3, 'User masuk ke Page PV Master List Maintenance', 'Pass', '2025-12-23T04:56:21.238Z'),
    @(4, 'User ID Input', 'Pass', '2025-12-23T04:56:21.867Z'),
    @(5, 'Password has been inputed', 'Pass', '2025-12-23T04:56:22.102Z'),
    @(6, 'Label Company is shown', 'Pass', '2025-12-23T04:56:35.863Z'),
    @(7, 'Company filter option exists: All', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(8, 'Company filter option exists: PT. BALI TELEKOM', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(9, 'Company filter option exists: PT. BATAVIA TOWERINDO', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(10, 'Company filter option exists: PT. INFRASTRUCTURE DIGITAL INDONESIA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(11, 'Company filter option exists: PT. JARINGAN PINTAR INDONESIA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(12, 'Company filter option exists: PT. MENARA BERSAMA TERPADU', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(13, 'Company filter option exists: PT. METRIC SOLUSI INTEGRASI', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(14, 'Company filter option exists: PT. MITRAYASA SARANA INFORMASI', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(15, 'Company filter option exists: PT. PERMATA KARYA PERDANA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(16, 'Company filter option exists: PT. PRIMA MEDIA SELARAS', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(17, 'Company filter option exists: PT. SOLUSINDO KREASI PRATAMA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(18, 'Company filter option exists: PT. SOLUSI MENARA INDONESIA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(19, 'Company filter option exists: PT. TOWER BERSAMA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(20, 'Company filter option exists: PT. TOWER BERSAMA INFRASTRUCTURE', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(21, 'Company filter option exists: PT. TELENET INTERNUSA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(22, 'Company filter option exists: PT. TOWERINDO KONVERGENSI', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(23, 'Company filter option exists: PT. TOWER ONE', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(24, 'Company filter option exists: PT. TRIAKA BERSAMA', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(25, 'Company filter option exists: PT. UNITED TOWERINDO', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(26, 'Company filter option exists: PT. VISI TELEKOMUNIKASI INFTRASTRUKTUR', 'Pass', '2025-12-23T04:56:35.904Z'),
    @(27, 'Button Clicked', 'Pass', '2025-12-23T04:56:53.523Z'),
    @(28, 'User masuk ke Page PV Master List Maintenance', 'Pass', '2025-12-23T04:56:53.523Z'),
    @(29, 'User ID Input', 'Pass', '2025-12-23T04:56:53.915Z'),
    @(30, 'Password has been inputed', 'Pass', '2025-12-23T04:56:54.146Z'),
    @(31, 'Radio "No" can be clicked via label', 'Pass', '2025-12-23T04:57:10.300Z'),
    @(32, 'Search button can be clicked', 'Pass', '2025-12-23T04:57:10.640Z'),
    @(33, 'Maintenance table loading spinner disappears', 'Pass', '2025-12-23T04:57:13.318Z'),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
